# Update Equipment Load Factors values (RMI files through Dec 1)

$wb = $excel.ActiveWorkbook

# --- ELF-bldg-winter ---
$wsWinter = $wb.Worksheets.Item("ELF-bldg-winter")
$wsWinter.Range("B2").Value = 4.9836
$wsWinter.Range("D2").Value = 4.42685
$wsWinter.Range("B5").Value = 1.20393
$wsWinter.Range("D5").Value = 1.11134
$wsWinter.Range("D7").Value = 1.11134

# --- ELF-bldg-summer ---
$wsSummer = $wb.Worksheets.Item("ELF-bldg-summer")
$wsSummer.Range("B3").Value = 7.22244
$wsSummer.Range("D3").Value = 7.38051
$wsSummer.Range("B5").Value = 1.38062
$wsSummer.Range("D5").Value = 1.59906
$wsSummer.Range("D7").Value = 1.59906

# --- ELF-vehicles ---
$wsVehicles = $wb.Worksheets.Item("ELF-vehicles")
$wsVehicles.Range("B4").Value = 0.93833
$wsVehicles.Range("C4").Value = 0.7369
$wsVehicles.Range("B5").Value = 0.93833
$wsVehicles.Range("C5").Value = 0.7369
$wsVehicles.Range("B6").Value = 0.93833
$wsVehicles.Range("C6").Value = 0.7369
$wsVehicles.Range("B7").Value = 0.93833
$wsVehicles.Range("C7").Value = 0.7369
